# Auto-generated Excel COM-interop script to apply profit-sheet updates
# across the Mateus_Profits workbook tables (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# --- Sheet 1: ALC ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("H141").Value = 5593.1
$ws.Range("I141").Value = 4818.4116
$ws.Range("K141").Value = 14455.2348
$ws.Range("M141").Value = -9275.234800000002

# --- Sheet 2: ARM ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("H32").Value = 1458.1698
$ws.Range("I32").Value = 1123.196
$ws.Range("K32").Value = 1123.196
$ws.Range("M32").Value = -836.1959999999999
$ws.Range("H50").Value = 22048
$ws.Range("I50").Value = 8074
$ws.Range("K50").Value = 8074
$ws.Range("M50").Value = -7360
$ws.Range("H61").Value = 4768466.5
$ws.Range("I61").Value = 7251319.5
$ws.Range("K61").Value = 7251319.5
$ws.Range("M61").Value = -7251107.5
$ws.Range("H74").Value = 3554.75
$ws.Range("I74").Value = 2674.3462
$ws.Range("K74").Value = 2674.3462
$ws.Range("M74").Value = -1800.3462
$ws.Range("H77").Value = 3554.75
$ws.Range("I77").Value = 2674.3462
$ws.Range("K77").Value = 13371.731
$ws.Range("M77").Value = -9003.731
$ws.Range("H88").Value = 1656.8572
$ws.Range("I88").Value = 1431.3334
$ws.Range("J88").Value = 1826
$ws.Range("K88").Value = 1431.3334
$ws.Range("L88").Value = 1826
$ws.Range("M88").Value = -1025.3334
$ws.Range("N88").Value = -2638
$ws.Range("H91").Value = 1656.8572
$ws.Range("I91").Value = 1431.3334
$ws.Range("J91").Value = 1826
$ws.Range("K91").Value = 1431.3334
$ws.Range("L91").Value = 1826
$ws.Range("M91").Value = -27.33339999999998
$ws.Range("N91").Value = -4634
$ws.Range("H122").Value = 3920.6
$ws.Range("I122").Value = 3920.6
$ws.Range("K122").Value = 11761.8
$ws.Range("M122").Value = -9311.799999999999
$ws.Range("H136").Value = 4768466.5
$ws.Range("I136").Value = 7251319.5
$ws.Range("K136").Value = 21753958.5
$ws.Range("M136").Value = -21751408.5

# --- Sheet 3: BSM ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("H20").Value = 2140.5454
$ws.Range("I20").Value = 2271.7778
$ws.Range("K20").Value = 2271.7778
$ws.Range("M20").Value = -2024.7778
$ws.Range("H86").Value = 11766329
$ws.Range("I86").Value = 1707.3846
$ws.Range("K86").Value = 1707.3846
$ws.Range("M86").Value = -584.3846000000001
$ws.Range("H89").Value = 11766329
$ws.Range("I89").Value = 1707.3846
$ws.Range("K89").Value = 8536.923000000001
$ws.Range("M89").Value = -2920.923000000001
$ws.Range("H99").Value = 3780.2195
$ws.Range("I99").Value = 2583.5
$ws.Range("K99").Value = 2583.5
$ws.Range("M99").Value = -1085.5
$ws.Range("H134").Value = 3739.2703
$ws.Range("I134").Value = 3787.5833
$ws.Range("K134").Value = 11362.7499
$ws.Range("M134").Value = -8827.749899999999

# --- Sheet 4: CRP ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("H4").Value = 10000
$ws.Range("J4").Value = 10000
$ws.Range("L4").Value = 10000
$ws.Range("N4").Value = -10224
$ws.Range("H68").Value = 45996.668
$ws.Range("J68").Value = 45996.668
$ws.Range("L68").Value = 45996.668
$ws.Range("N68").Value = -47494.668
$ws.Range("H71").Value = 45996.668
$ws.Range("J71").Value = 45996.668
$ws.Range("L71").Value = 137990.004
$ws.Range("N71").Value = -145478.004
$ws.Range("H141").Value = 37699.5
$ws.Range("I141").Value = 37699
$ws.Range("J141").Value = 37700
$ws.Range("K141").Value = 37699
$ws.Range("L141").Value = 37700
$ws.Range("M141").Value = -32519
$ws.Range("N141").Value = -48060

# --- Sheet 5: CUL ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("H131").Value = 20001554
$ws.Range("J131").Value = 6252532
$ws.Range("L131").Value = 18757596
$ws.Range("N131").Value = -18767676

# --- Sheet 6: GSM ---
$ws = $wb.Worksheets.Item(6)
$ws.Range("H70").Value = 17999.857
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 17999.857
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 17999.857
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -18539.857
$ws.Range("H73").Value = 17999.857
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 17999.857
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 17999.857
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -19871.857
$ws.Range("H113").Value = 1000995
$ws.Range("I113").Value = 1000995
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1000995
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -998825
$ws.Range("N113").ClearContents()
$ws.Range("H122").Value = 2274.875
$ws.Range("I122").Value = 2274.875
$ws.Range("K122").Value = 6824.625
$ws.Range("M122").Value = -4374.625
$ws.Range("H126").Value = 6344.4
$ws.Range("I126").Value = 7375
$ws.Range("K126").Value = 22125
$ws.Range("M126").Value = -19655
$ws.Range("H132").Value = 5800.2915
$ws.Range("I132").Value = 5464
$ws.Range("J132").Value = 9499.5
$ws.Range("K132").Value = 16392
$ws.Range("L132").Value = 28498.5
$ws.Range("M132").Value = -13862
$ws.Range("N132").Value = -33558.5

# --- Sheet 7: LTW ---
$ws = $wb.Worksheets.Item(7)
$ws.Range("H61").Value = 253252
$ws.Range("I61").Value = 253252
$ws.Range("K61").Value = 253252
$ws.Range("M61").Value = -253050
$ws.Range("H68").Value = 8281.5625
$ws.Range("J68").Value = 7460.875
$ws.Range("L68").Value = 7460.875
$ws.Range("N68").Value = -8958.875
$ws.Range("H71").Value = 8281.5625
$ws.Range("J71").Value = 7460.875
$ws.Range("L71").Value = 37304.375
$ws.Range("N71").Value = -44792.375
$ws.Range("H82").Value = 4630.2085
$ws.Range("I82").Value = 3494
$ws.Range("J82").Value = 5591.615
$ws.Range("K82").Value = 3494
$ws.Range("L82").Value = 5591.615
$ws.Range("M82").Value = -3133
$ws.Range("N82").Value = -6313.615
$ws.Range("H85").Value = 4630.2085
$ws.Range("I85").Value = 3494
$ws.Range("J85").Value = 5591.615
$ws.Range("K85").Value = 3494
$ws.Range("L85").Value = 5591.615
$ws.Range("M85").Value = -2246
$ws.Range("N85").Value = -8087.615
$ws.Range("H93").Value = 16488.5
$ws.Range("I93").Value = 985.63635
$ws.Range("K93").Value = 985.63635
$ws.Range("M93").Value = 262.36365
$ws.Range("H100").Value = 2276238
$ws.Range("I100").Value = 3127952.5
$ws.Range("K100").Value = 3127952.5
$ws.Range("M100").Value = -3127411.5
$ws.Range("H113").Value = 253252
$ws.Range("I113").Value = 253252
$ws.Range("K113").Value = 253252
$ws.Range("M113").Value = -251082
$ws.Range("H132").Value = 10242.5625
$ws.Range("I132").Value = 11782
$ws.Range("J132").Value = 6855.8
$ws.Range("K132").Value = 35346
$ws.Range("L132").Value = 20567.4
$ws.Range("M132").Value = -32816
$ws.Range("N132").Value = -25627.4

# --- Sheet 8: WVR ---
$ws = $wb.Worksheets.Item(8)
$ws.Range("H132").Value = 5088.844
$ws.Range("I132").Value = 4572.484
$ws.Range("J132").Value = 6232.2144
$ws.Range("K132").Value = 13717.452
$ws.Range("L132").Value = 18696.6432
$ws.Range("M132").Value = -11187.452
$ws.Range("N132").Value = -23756.6432

Write-Host "Applied all profit-sheet updates."